# Append closed/open trade #41 (leadlag UP) as new row 31 on the
# "leadlag" worksheet, growing the used range from A1:N30 to A1:N31.
#
# NOTE: Date/empty-string fields are prefixed with a leading apostrophe.
# Without it, Excel's COM Value setter auto-converts a "YYYY-MM-DD"
# string into a real date serial, and assigning "" simply clears/removes
# the cell instead of leaving a present-but-empty text cell. The leading
# apostrophe forces literal text entry, matching the source data's text
# columns (Date/Time/Strategy/Side/Status/Entry Reason/Exit Reason are
# all stored as text, not numbers/dates, in this workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 31

$ws.Cells.Item($row, 1).Value  = 41                                         # Trade #
$ws.Cells.Item($row, 2).Value  = "'2026-02-16"                              # Date (text)
$ws.Cells.Item($row, 3).Value  = "21:29:18"                                 # Time
$ws.Cells.Item($row, 4).Value  = "leadlag"                                  # Strategy
$ws.Cells.Item($row, 5).Value  = "UP"                                       # Side
$ws.Cells.Item($row, 6).Value  = 68714.66                                   # Entry Price
$ws.Cells.Item($row, 7).Value  = "'"                                        # Exit Price (blank text)
$ws.Cells.Item($row, 8).Value  = "OPEN"                                     # Status
$ws.Cells.Item($row, 9).Value  = 0                                          # P&L %
$ws.Cells.Item($row, 10).Value = 0                                          # P&L $
$ws.Cells.Item($row, 11).Value = 0.75                                       # Confidence
$ws.Cells.Item($row, 12).Value = "Binance leading with 0.114% move"         # Entry Reason
$ws.Cells.Item($row, 13).Value = "'"                                        # Exit Reason (blank text)
$ws.Cells.Item($row, 14).Value = 0                                          # Duration (min)
